# Applies edits described in the commit:
# "rajout titre page accueil et keywords page2"
#
# 1. Fix a typo in an existing string: "Manifier" -> "Magnifier"
#    ("Magnifier et compresser les ressources (images et code)")
# 2. Add a new row (row 20) to the audit table with:
#      Catégorie = Accessibilité
#      Problème identifié = Langue
#      Explication du problème = Les documents html ont comme paramètre de langue "Default"
#      Bonne pratique/Action recommandée = Changer ce paramètre pour mettre "fr" à la place
#      Action recommandée / résolu checkbox = TRUE
# 3. Move the selected/active cell to A15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new row 20 content ---
$ws.Range("A20").Value = "Accessibilité"
$ws.Range("B20").Value = "Langue"
$ws.Range("C20").Value = "Les documents html ont comme paramètre de langue ""Default"""
$ws.Range("D20").Value = "Changer ce paramètre pour mettre ""fr"" à la place"
$ws.Range("E20").Value = $true

# Copy formatting from the row above (row 19) so the new row matches
# the rest of the table's look and feel.
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122) # xlPasteFormats

# --- 2. Fix typo in existing cell B9 ---
$ws.Range("B9").Value = "Magnifier et compresser les ressources (images et code)"

# --- 3. Update the active selection ---
$ws.Range("A15").Select()
